$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19, 8).Value = 7210
$ws.Cells.Item(19, 9).Value = 3683.3333
$ws.Cells.Item(19, 11).Value = 3683.3333
$ws.Cells.Item(19, 13).Value = -3508.3333
# Row 132
$ws.Cells.Item(132, 8).Value = 3322.8572
$ws.Cells.Item(132, 9).Value = 2646.7805
$ws.Cells.Item(132, 10).Value = 6787.75
$ws.Cells.Item(132, 11).Value = 7940.3415
$ws.Cells.Item(132, 12).Value = 20363.25
$ws.Cells.Item(132, 13).Value = -5410.3415
$ws.Cells.Item(132, 14).Value = -25423.25
# Row 138
$ws.Cells.Item(138, 8).Value = 2279.3704
$ws.Cells.Item(138, 9).Value = 1244.7
$ws.Cells.Item(138, 10).Value = 3572.7083
$ws.Cells.Item(138, 11).Value = 3734.1
$ws.Cells.Item(138, 12).Value = 10718.1249
$ws.Cells.Item(138, 13).Value = 1405.9
$ws.Cells.Item(138, 14).Value = -20998.1249
# Row 141
$ws.Cells.Item(141, 8).Value = 3523.75
$ws.Cells.Item(141, 9).Value = 3063.3333
$ws.Cells.Item(141, 10).Value = 4905
$ws.Cells.Item(141, 11).Value = 9189.999899999999
$ws.Cells.Item(141, 12).Value = 14715
$ws.Cells.Item(141, 13).Value = -4009.999899999999
$ws.Cells.Item(141, 14).Value = -25075

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 7062.66
$ws.Cells.Item(32, 9).Value = 3505.64
$ws.Cells.Item(32, 10).Value = 17733.72
$ws.Cells.Item(32, 11).Value = 3505.64
$ws.Cells.Item(32, 12).Value = 17733.72
$ws.Cells.Item(32, 13).Value = -3218.64
$ws.Cells.Item(32, 14).Value = -18307.72
# Row 61
$ws.Cells.Item(61, 8).Value = 8774259
$ws.Cells.Item(61, 9).Value = 11113535
$ws.Cells.Item(61, 10).Value = 1974.75
$ws.Cells.Item(61, 11).Value = 11113535
$ws.Cells.Item(61, 12).Value = 1974.75
$ws.Cells.Item(61, 13).Value = -11113323
$ws.Cells.Item(61, 14).Value = -2398.75
# Row 111
$ws.Cells.Item(111, 8).Value = 40000
$ws.Cells.Item(111, 10).Value = 40000
$ws.Cells.Item(111, 12).Value = 40000
$ws.Cells.Item(111, 14).Value = -48180
# Row 117
$ws.Cells.Item(117, 8).Value = 30171.77
$ws.Cells.Item(117, 10).Value = 30171.77
$ws.Cells.Item(117, 12).Value = 30171.77
$ws.Cells.Item(117, 14).Value = -39349.77
# Row 123
$ws.Cells.Item(123, 8).Value = 31666.666
$ws.Cells.Item(123, 10).Value = 31666.666
$ws.Cells.Item(123, 12).Value = 31666.666
$ws.Cells.Item(123, 14).Value = -41466.666
# Row 136
$ws.Cells.Item(136, 8).Value = 8774259
$ws.Cells.Item(136, 9).Value = 11113535
$ws.Cells.Item(136, 10).Value = 1974.75
$ws.Cells.Item(136, 11).Value = 33340605
$ws.Cells.Item(136, 12).Value = 5924.25
$ws.Cells.Item(136, 13).Value = -33338055
$ws.Cells.Item(136, 14).Value = -11024.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 118
$ws.Cells.Item(118, 8).Value = 7870
$ws.Cells.Item(118, 10).Value = 7870
$ws.Cells.Item(118, 12).Value = 7870
$ws.Cells.Item(118, 14).Value = -11184
# Row 134
$ws.Cells.Item(134, 8).Value = 4817.6665
$ws.Cells.Item(134, 9).Value = 4195.7856
$ws.Cells.Item(134, 11).Value = 12587.3568
$ws.Cells.Item(134, 13).Value = -10052.3568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 1890.6
$ws.Cells.Item(16, 9).Value = 1877.25
$ws.Cells.Item(16, 10).Value = 1899.5
$ws.Cells.Item(16, 11).Value = 1877.25
$ws.Cells.Item(16, 12).Value = 1899.5
$ws.Cells.Item(16, 13).Value = -1590.25
$ws.Cells.Item(16, 14).Value = -2473.5
# Row 31
$ws.Cells.Item(31, 8).Value = 4507453
$ws.Cells.Item(31, 9).Value = 2832.1077
$ws.Cells.Item(31, 11).Value = 2832.1077
$ws.Cells.Item(31, 13).Value = -2537.1077
# Row 34
$ws.Cells.Item(34, 8).Value = 4507453
$ws.Cells.Item(34, 9).Value = 2832.1077
$ws.Cells.Item(34, 11).Value = 2832.1077
$ws.Cells.Item(34, 13).Value = -2630.1077
# Row 58
$ws.Cells.Item(58, 8).Value = 2107.5483
$ws.Cells.Item(58, 9).Value = 1233.2222
$ws.Cells.Item(58, 10).Value = 3318.1538
$ws.Cells.Item(58, 11).Value = 1233.2222
$ws.Cells.Item(58, 12).Value = 3318.1538
$ws.Cells.Item(58, 13).Value = -1030.2222
$ws.Cells.Item(58, 14).Value = -3724.1538
# Row 99
$ws.Cells.Item(99, 8).Value = 1619.9565
$ws.Cells.Item(99, 9).Value = 1428.6842
$ws.Cells.Item(99, 10).Value = 2528.5
$ws.Cells.Item(99, 11).Value = 1428.6842
$ws.Cells.Item(99, 12).Value = 2528.5
$ws.Cells.Item(99, 13).Value = 69.31580000000008
$ws.Cells.Item(99, 14).Value = -5524.5
# Row 110
$ws.Cells.Item(110, 8).Value = 37701.777
$ws.Cells.Item(110, 10).Value = 37701.777
$ws.Cells.Item(110, 12).Value = 37701.777
$ws.Cells.Item(110, 14).Value = -45881.777
# Row 113
$ws.Cells.Item(113, 8).Value = 1890.6
$ws.Cells.Item(113, 9).Value = 1877.25
$ws.Cells.Item(113, 10).Value = 1899.5
$ws.Cells.Item(113, 11).Value = 1877.25
$ws.Cells.Item(113, 12).Value = 1899.5
$ws.Cells.Item(113, 13).Value = 292.75
$ws.Cells.Item(113, 14).Value = -6239.5
# Row 126
$ws.Cells.Item(126, 8).Value = 1619.9565
$ws.Cells.Item(126, 9).Value = 1428.6842
$ws.Cells.Item(126, 10).Value = 2528.5
$ws.Cells.Item(126, 11).Value = 4286.0526
$ws.Cells.Item(126, 12).Value = 7585.5
$ws.Cells.Item(126, 13).Value = -1816.0526
$ws.Cells.Item(126, 14).Value = -12525.5
# Row 136
$ws.Cells.Item(136, 8).Value = 2107.5483
$ws.Cells.Item(136, 9).Value = 1233.2222
$ws.Cells.Item(136, 10).Value = 3318.1538
$ws.Cells.Item(136, 11).Value = 3699.6666
$ws.Cells.Item(136, 12).Value = 9954.4614
$ws.Cells.Item(136, 13).Value = -1149.6666
$ws.Cells.Item(136, 14).Value = -15054.4614

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 3851.6667
$ws.Cells.Item(102, 9).Value = 5710.5454
$ws.Cells.Item(102, 10).Value = 930.5714
$ws.Cells.Item(102, 11).Value = 5710.5454
$ws.Cells.Item(102, 12).Value = 930.5714
$ws.Cells.Item(102, 13).Value = -4088.5454
$ws.Cells.Item(102, 14).Value = -4174.5714
# Row 118
$ws.Cells.Item(118, 8).Value = 17516.666
$ws.Cells.Item(118, 10).Value = 17516.666
$ws.Cells.Item(118, 12).Value = 17516.666
$ws.Cells.Item(118, 14).Value = -20830.666
# Row 122
$ws.Cells.Item(122, 8).Value = 4447767.5
$ws.Cells.Item(122, 9).Value = 11113903
$ws.Cells.Item(122, 10).Value = 3677.3333
$ws.Cells.Item(122, 11).Value = 33341709
$ws.Cells.Item(122, 12).Value = 11031.9999
$ws.Cells.Item(122, 13).Value = -33339259
$ws.Cells.Item(122, 14).Value = -15931.9999
# Row 126
$ws.Cells.Item(126, 8).Value = 3411
$ws.Cells.Item(126, 9).Value = 1452.5
$ws.Cells.Item(126, 10).Value = 4716.6665
$ws.Cells.Item(126, 11).Value = 4357.5
$ws.Cells.Item(126, 12).Value = 14149.9995
$ws.Cells.Item(126, 13).Value = -1887.5
$ws.Cells.Item(126, 14).Value = -19089.9995
# Row 132
$ws.Cells.Item(132, 8).Value = 4095.2
$ws.Cells.Item(132, 9).Value = 5259.9033
$ws.Cells.Item(132, 10).Value = 2194.8948
$ws.Cells.Item(132, 11).Value = 15779.7099
$ws.Cells.Item(132, 12).Value = 6584.6844
$ws.Cells.Item(132, 13).Value = -13249.7099
$ws.Cells.Item(132, 14).Value = -11644.6844

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 7229
$ws.Cells.Item(7, 9).Value = 14750
$ws.Cells.Item(7, 11).Value = 14750
$ws.Cells.Item(7, 13).Value = -14638
# Row 108
$ws.Cells.Item(108, 8).Value = 25989.5
$ws.Cells.Item(108, 10).Value = 25989.5
$ws.Cells.Item(108, 12).Value = 25989.5
$ws.Cells.Item(108, 14).Value = -33669.5
# Row 122
$ws.Cells.Item(122, 8).Value = 6823.727
$ws.Cells.Item(122, 9).Value = 9326.5
$ws.Cells.Item(122, 10).Value = 5393.5713
$ws.Cells.Item(122, 11).Value = 27979.5
$ws.Cells.Item(122, 12).Value = 16180.7139
$ws.Cells.Item(122, 13).Value = -25529.5
$ws.Cells.Item(122, 14).Value = -21080.7139
# Row 126
$ws.Cells.Item(126, 8).Value = 7229
$ws.Cells.Item(126, 9).Value = 14750
$ws.Cells.Item(126, 11).Value = 44250
$ws.Cells.Item(126, 13).Value = -41780
# Row 132
$ws.Cells.Item(132, 8).Value = 11369797
$ws.Cells.Item(132, 9).Value = 2802.6206
$ws.Cells.Item(132, 10).Value = 33345986
$ws.Cells.Item(132, 11).Value = 8407.861800000001
$ws.Cells.Item(132, 12).Value = 100037958
$ws.Cells.Item(132, 13).Value = -5877.861800000001
$ws.Cells.Item(132, 14).Value = -100043018

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Cells.Item(45, 8).Value = 6626
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 6626
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 6626
$ws.Cells.Item(45, 14).Value = -7608
$ws.Cells.Item(45, 13).ClearContents()
# Row 122
$ws.Cells.Item(122, 8).Value = 2532.875
$ws.Cells.Item(122, 9).Value = 2606.5
$ws.Cells.Item(122, 10).Value = 2385.625
$ws.Cells.Item(122, 11).Value = 7819.5
$ws.Cells.Item(122, 12).Value = 7156.875
$ws.Cells.Item(122, 13).Value = -5369.5
$ws.Cells.Item(122, 14).Value = -12056.875
# Row 126
$ws.Cells.Item(126, 8).Value = 1822.2858
$ws.Cells.Item(126, 9).Value = 1531.2
$ws.Cells.Item(126, 10).Value = 2550
$ws.Cells.Item(126, 11).Value = 4593.6
$ws.Cells.Item(126, 12).Value = 7650
$ws.Cells.Item(126, 13).Value = -2123.6
$ws.Cells.Item(126, 14).Value = -12590
